# Correct closing-date values for the algo (commit: "closing dates were wrong for algo")
$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("M&MFIN")

# Row 7
$ws.Range("F7").Value = 267.95
$ws.Range("G7").Value = 272.55
$ws.Range("H7").Value = 266.2
$ws.Range("I7").Value = 269.9
$ws.Range("J7").Value = 268.8

# Row 9
$ws.Range("G9").Value = 269
$ws.Range("H9").Value = 261.9
$ws.Range("I9").Value = 267.2

# Row 10
$ws.Range("G10").Value = 270.75
$ws.Range("H10").Value = 266.3
$ws.Range("I10").Value = 270.2

# Row 11
$ws.Range("G11").Value = 271.95
$ws.Range("H11").Value = 269.75
$ws.Range("I11").Value = 270.9

# Row 12
$ws.Range("G12").Value = 272.3
$ws.Range("H12").Value = 269.5
$ws.Range("I12").Value = 271.3

# Row 13
$ws.Range("G13").Value = 271.7
$ws.Range("H13").Value = 270.1
$ws.Range("I13").Value = 270.1

# Row 14
$ws.Range("G14").Value = 271.15
$ws.Range("H14").Value = 270
$ws.Range("I14").Value = 270.8

# Row 15
$ws.Range("G15").Value = 271.65
$ws.Range("H15").Value = 270.4
$ws.Range("I15").Value = 271

# Row 16
$ws.Range("G16").Value = 272.4
$ws.Range("H16").Value = 270.85
$ws.Range("I16").Value = 272.25

# Row 17
$ws.Range("G17").Value = 272.35
$ws.Range("H17").Value = 271.05
$ws.Range("I17").Value = 272.3

# Row 18
$ws.Range("G18").Value = 272.55
$ws.Range("H18").Value = 270
$ws.Range("I18").Value = 270.25

# Row 19
$ws.Range("G19").Value = 270.65
$ws.Range("H19").Value = 269.25
$ws.Range("I19").Value = 270.05

# Row 20
$ws.Range("G20").Value = 270.35
$ws.Range("H20").Value = 268.2
$ws.Range("I20").Value = 269.8

# Row 21
$ws.Range("G21").Value = 270.95
$ws.Range("H21").Value = 269
$ws.Range("I21").Value = 269.95
